$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking strings so Excel does not coerce them to floats
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D25","D27","D29","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "28.619.83"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "1.850.77"
$ws.Range("E3").Value = "  -3.93%  "
$ws.Range("E4").Value = "  -1.10%  "
$ws.Range("D5").Value = "335.61"
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "0.4666"
$ws.Range("E7").Value = "  -3.15%  "
$ws.Range("E8").Value = "  -3.75%  "
$ws.Range("D9").Value = "46.25"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("D10").Value = "0.07903"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("D11").Value = "0.9779"
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("D12").Value = "22.29"
$ws.Range("E12").Value = "  -6.31%  "
$ws.Range("D13").Value = "1.884.28"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "5.817"
$ws.Range("E14").Value = "  -4.28%  "
$ws.Range("D15").Value = "6.963"
$ws.Range("E15").Value = "  -4.52%  "
$ws.Range("D16").Value = "0.06908"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "87.79"
$ws.Range("E17").Value = "  -4.05%  "
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "0.00001003"
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "28.637.66"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("E23").Value = "  -4.69%  "
$ws.Range("E24").Value = "  -6.12%  "
$ws.Range("D25").Value = "2.158"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "2.119.87"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "153.16"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("D29").Value = "6.055"
$ws.Range("E29").Value = "  -4.76%  "
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("D31").Value = "117.53"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("D32").Value = "0.9677"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("D33").Value = "0.09347"
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("D34").Value = "5.363"
$ws.Range("E34").Value = "  -4.53%  "
$ws.Range("D35").Value = "3.476"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("D37").Value = "0.06105"
$ws.Range("E37").Value = "  -6.62%  "
$ws.Range("D38").Value = "0.02200"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("D39").Value = "1.165"
$ws.Range("E39").Value = "  -4.20%  "
$ws.Range("D40").Value = "7.683"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").Value = "0.5693"
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("D42").Value = "10.09"
$ws.Range("E42").Value = "  -5.98%  "
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").Value = "2.422"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "1.246"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5372"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "11.73"
$ws.Range("E47").Value = "  -4.92%  "
$ws.Range("D48").Value = "0.07094"
$ws.Range("E48").Value = "  -5.96%  "
$ws.Range("D49").Value = "1.904"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("D50").Value = "113.17"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").Value = "2.343"
$ws.Range("E51").Value = "  -3.58%  "

# Reset style index back to default (Normal) for cells where we forced text format,
# so no spurious style id is introduced, matching original formatting
foreach ($addr in $textCells) { $ws.Range($addr).Style = "Normal" }
